$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.783.25'
$ws.Range("E2").Value = '  -3.94%  '
$ws.Range("D3").Value = '1.817.40'
$ws.Range("E3").Value = '  -2.86%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '277.65'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -7.82%  '
$ws.Range("E6").Value = '  -0.05%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5095'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -4.40%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3519'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -6.15%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '45.41'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.23%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.06657'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -7.20%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '19.96'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -7.32%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.8333'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -6.03%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.07911'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -3.02%  '
$ws.Range("D14").Value = '1.831.23'
$ws.Range("E14").Value = '  -2.19%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.080'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -3.28%  '
$ws.Range("E16").Value = '  -5.78%  '
$ws.Range("E17").Value = '  -0.10%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '14.10'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -3.97%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000008024'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -6.00%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.000'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.07%  '
$ws.Range("D21").Value = '25.834.15'
$ws.Range("E21").Value = '  -3.93%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.728'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -4.86%  '
$ws.Range("E23").Value = '  -6.39%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.084'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -4.79%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '142.37'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -3.18%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.189'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -3.15%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.668'
$ws.Range("D27").Style = "Normal"
$ws.Range("E28").Value = '  -5.02%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '109.62'
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.342'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -8.47%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.240'
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.08833'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -3.27%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.04867'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -2.66%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7331'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -8.20%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.132'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -3.22%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.880'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -3.61%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.163'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.48%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.9998'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.03%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.5225'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -13.16%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.332'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -10.00%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.01850'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -5.16%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.9549'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -10.93%  '
$ws.Range("B43").Value = 'FraxShare'
$ws.Range("C43").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.214'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -6.37%  '
$ws.Range("B44").Value = 'Quant'
$ws.Range("C44").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '111.47'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -3.76%  '
$ws.Range("E45").Value = '  -9.21%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.000'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.03%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.4610'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -10.11%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.1365'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -8.74%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '36.56'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -2.91%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '9.180'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -7.28%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.502'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -7.61%  '
